# Generate Report for Handoff
#
# For the rows that are "Ready for handoff" (39c51c0a, 78cb526c, 7b72c2d5,
# 7f0472fa, e5c59853, ee176d73 -> rows 7,8,9,10,13,14 on every sheet) the
# handoff report generation now records:
#   - a Priority of "ht" (previously blank) on the zh-cn and de-de sheets
#   - a refreshed "Latest Handoff Datetime" / "Latest HO Xliff Generate Date"
#     timestamp on every sheet that shows it

$wb = $excel.ActiveWorkbook

$overview = $wb.Worksheets.Item("Overview")
$zhcn = $wb.Worksheets.Item("zh-cn")
$dede = $wb.Worksheets.Item("de-de")

$rows = @(7, 8, 9, 10, 13, 14)

foreach ($r in $rows) {
    # Overview sheet: column G = "Latest HO Xliff Generate Date"
    $overview.Range("G$r").Value = "2016-08-25 08:22:34"

    # zh-cn sheet: column E = "Priority", column H = "Latest Handoff Datetime"
    $zhcn.Range("E$r").Value = "ht"
    $zhcn.Range("H$r").Value = "2016-08-25 08:22:29"

    # de-de sheet: column E = "Priority", column H = "Latest Handoff Datetime"
    $dede.Range("E$r").Value = "ht"
    $dede.Range("H$r").Value = "2016-08-25 08:22:34"
}
